# Add a "Save" column (H) to the s_vals sheet, matching the header style
# already used by the other header cells (copy G1's format onto H1), and
# add the corresponding data value in H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header cell's formatting (bold, bordered, centered) onto H1.
$ws.Range("G1").Copy($ws.Range("H1"))

# Set the new header text and the data value for row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
